$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.210.95'
$ws.Range("E2").Value = '  -5.96%  '

$ws.Range("D3").Value = '2.558.49'
$ws.Range("E3").Value = '  -1.87%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '299.24'
$ws.Range("E5").Value = '  -3.28%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '93.11'
$ws.Range("E6").Value = '  -5.66%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.576'
$ws.Range("E7").Value = '  -3.45%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.553'
$ws.Range("E9").Value = '  -4.47%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.05'
$ws.Range("E10").Value = '  -6.98%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0812'
$ws.Range("E11").Value = '  -3.59%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.79'
$ws.Range("E12").Value = '  -3.48%  '

$ws.Range("E13").Value = '  +1.72%  '

$ws.Range("D14").Value = '2.941.05'
$ws.Range("E14").Value = '  -2.45%  '

$ws.Range("D15").Value = '2.544.26'
$ws.Range("E15").Value = '  -2.68%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.876'
$ws.Range("E16").Value = '  -3.84%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.23'
$ws.Range("E17").Value = '  -3.78%  '

$ws.Range("D18").Value = '43.226.71'
$ws.Range("E18").Value = '  -6.46%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.21'
$ws.Range("E19").Value = '  +4.02%  '

$ws.Range("D20").Value = '0.0₃0986'
$ws.Range("E20").Value = '  -2.75%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.67'
$ws.Range("E21").Value = '  -0.93%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '72.21'
$ws.Range("E22").Value = '  -1.04%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '261.46'
$ws.Range("E23").Value = '  -9.86%  '

$ws.Range("E24").Value = '  -4.09%  '

$ws.Range("B25").Value = 'EthereumClassic'
$ws.Range("C25").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '29.67'
$ws.Range("E25").Value = '  +0.31%  '

$ws.Range("B26").Value = 'ImmutableX'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.15'
$ws.Range("E26").Value = '  -4.65%  '

$ws.Range("E27").Value = '  +0.20%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.09'
$ws.Range("E28").Value = '  -6.35%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '37.00'
$ws.Range("E29").Value = '  -4.43%  '

$ws.Range("E30").Value = '  -3.74%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.05'
$ws.Range("E31").Value = '  -3.17%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '155.96'
$ws.Range("E32").Value = '  -2.02%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.19'
$ws.Range("E33").Value = '  -2.31%  '

$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.40'
$ws.Range("E34").Value = '  -5.63%  '

$ws.Range("B35").Value = 'WEMIXToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.75'
$ws.Range("E35").Value = '  -2.06%  '

$ws.Range("E36").Value = '  -4.81%  '

$ws.Range("E37").Value = '  -4.82%  '

$ws.Range("E38").Value = '  -2.60%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '16.73'
$ws.Range("E39").Value = '  +6.94%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '23.44'
$ws.Range("E40").Value = '  +8.88%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.51'
$ws.Range("E41").Value = '  -1.03%  '

$ws.Range("E42").Value = '  -4.78%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.92'
$ws.Range("E43").Value = '  -2.37%  '

$ws.Range("D44").Value = '2.068.09'
$ws.Range("E44").Value = '  -2.65%  '

$ws.Range("E45").Value = '  -0.12%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '86.17'
$ws.Range("E46").Value = '  -10.58%  '

$ws.Range("E47").Value = '  +2.60%  '

$ws.Range("D48").Value = '2.798.46'
$ws.Range("E48").Value = '  -2.46%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.78'
$ws.Range("E49").Value = '  -7.47%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.71'
$ws.Range("E50").Value = '  -1.27%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '104.61'
$ws.Range("E51").Value = '  -4.59%  '
